# The team miscounted how many tasks they had: the sprint total goes
# from 96 to 101. That single source value drives the "Ideal" burndown
# line (column B) and the starting point of the "Actual" burndown line
# (column C). We also add the newly-tracked day (row 9) to the Actual
# column, since a task was logged on that day too.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sprint 4")
$ws.Activate()

# Correct total task count (was 96, should have been 101)
$ws.Range("B3").Value = 101
$ws.Range("C3").Value = 101

# A task-day that had been omitted now has its actual count tracked,
# following the same pattern as the cells above it.
$ws.Range("C9").Formula = "=C8-F9"

# Update the selection to match where the new row now sits.
$ws.Range("C8:C9").Select()
